# Daily attendance processing - rotate "Recorded By" (column G) author lists
# so the first (leading) author moves to the end of the comma-separated list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row based on column A (data starts at row 2; row 1 is the header)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) { continue }

    $parts = $val -split ", "

    if ($parts.Length -gt 1) {
        $rotated = $parts[1..($parts.Length - 1)] + $parts[0]
        $newVal = $rotated -join ", "
        $cell.Value = $newVal
    }
}
